# Applies the Domain_Objects.docx content edits described by the commit
# "hw creating entities / repositories":
#   - Product.product_type changes type from int -> str
#   - date(...) fields get a space before "(" and use "/" instead of "-"
#     as the day/month/year separator (dd/mm/yyyy instead of dd-mm-yyyy)
#   - Transaction.transaction_assets simplifies its type annotation from
#     list[Product.name, Product.quantity, Product.price] to list [Product]
#   - Invoice.invoicer_info is renamed to invoicer (and its "Counterparty"
#     type gets pushed two more tab-stops to line up with the other rows)
#   - Invoice.bill_to_info is renamed to invoice_to
#   - Invoice.items gets an explicit " (list [Product])" type annotation
#     appended after "Transactions.transaction_assets"

$d = $word.ActiveDocument
$tab = [char]9

function Replace-InParagraph {
    param(
        [int]$Index,
        [string]$FindText,
        [string]$ReplaceText,
        [bool]$MatchWholeWord
    )
    $p = $d.Paragraphs.Item($Index)
    $r = $p.Range
    $r.Find.Execute($FindText, $MatchWholeWord, $false, $false, $false, $false, $true, 1, $false, $ReplaceText, 2) | Out-Null
}

# --- Product: product_type  int -> str ---------------------------------
Replace-InParagraph 25 "int" "str" $true

# --- Transaction: transaction_date  date(dd-mm-yyyy ...) -> date (dd/mm/yyyy ...) ---
Replace-InParagraph 52 "date(" "date (" $false
Replace-InParagraph 52 "-mm-" "/mm/" $false

# --- Transaction: transaction_assets  list[...] -> list [Product] ------
Replace-InParagraph 55 "list[Product.name, Product.quantity, Product.price]" "list [Product]" $false

# --- Invoice: invoicer_info -> invoicer (+2 tabs before "Counterparty") ---
Replace-InParagraph 61 "invoicer_info" "invoicer" $true
Replace-InParagraph 61 "Counterparty" "$tab${tab}Counterparty" $false

# --- Invoice: bill_to_info -> invoice_to --------------------------------
Replace-InParagraph 62 "bill_to_info" "invoice_to" $true

# --- Invoice: invoice_date  date(dd-mm-yyyy ...) -> date (dd/mm/yyyy ...) ---
Replace-InParagraph 63 "date(" "date (" $false
Replace-InParagraph 63 "-mm-" "/mm/" $false

# --- Invoice: invoice_due_to  date(dd-mm-yyyy ...) -> date (dd/mm/yyyy ...) ---
Replace-InParagraph 64 "date(" "date (" $false
Replace-InParagraph 64 "-mm-" "/mm/" $false

# --- Invoice: items  append " (list [Product])" after Transactions.transaction_assets ---
Replace-InParagraph 65 "Transactions.transaction_assets" "Transactions.transaction_assets (list [Product])" $false
